# SBTM Report Template v1.1 - remove a stray/duplicated session-log entry
# from the "Summary" sheet (rows 21-23). Mirrors clearing the cell
# contents/cells in Excel for the extra "addGrade" session row that was
# duplicated into row 23, and the now-unneeded C/D values in rows 21-22.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# Row 21: clear the Sessions/Bugs tally values in C21/D21, but keep the
# (styled) blank cells in place.
$ws.Range("C21").ClearContents()
$ws.Range("D21").ClearContents()

# Row 22: C22 is removed outright (no longer referenced), D22 is cleared
# back to a styled blank.
$ws.Range("C22").Clear()
$ws.Range("D22").ClearContents()

# Row 23: this whole record (tester/charter/date/session length/bug count
# and the running-total formulas) was a duplicate of row 22 and is removed.
# H23/I23/J23/L23/M23 go away completely; K23/N23/O23 remain as empty,
# styled placeholder cells.
$ws.Range("H23").Clear()
$ws.Range("I23").Clear()
$ws.Range("J23").Clear()
$ws.Range("K23").ClearContents()
$ws.Range("L23").Clear()
$ws.Range("M23").Clear()
$ws.Range("N23").ClearContents()
$ws.Range("O23").ClearContents()
